$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Segunda (Monday) Metrologia 1 class slides over to other weekdays
$ws.Range("F2").Value = "[-, -, -, 'MEC-1B-Metrologia 1']"

$ws.Range("B3").Value = "-"
$ws.Range("E3").Value = "[-, -, 'MEC-1B-Metrologia 1', -]"

$ws.Range("B4").Value = "-"
$ws.Range("E4").Value = "[-, -, 'MEC-1B-Metrologia 1', -]"

$ws.Range("B6").Value = "-"

$ws.Range("B7").Value = "-"

$ws.Range("F8").Value = "['MEC-1B-Metrologia 1', -, -, -]"

# Night classes (rows 18-21)
$ws.Range("B18").Value = "['ELM-1NA-Metrologia', -, -, -]"
$ws.Range("C18").Value = "[-, -, -, 'MEC-1NB-Metrologia 1']"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "['ELM-1NA-Metrologia', -, -, -]"

$ws.Range("B19").Value = "['ELM-1NA-Metrologia', -, -, -]"
$ws.Range("C19").Value = "[-, -, -, 'MEC-1NB-Metrologia 1']"
$ws.Range("E19").Value = "['ELM-1NA-Metrologia', -, -, -]"

$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "[-, -, -, 'MEC-1NB-Metrologia 1']"
$ws.Range("E20").Value = "['ELM-2NA-CAD', 'ELM-2NA-CAD']"
$ws.Range("F20").Value = "-"

$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "[-, -, -, 'MEC-1NB-Metrologia 1']"
$ws.Range("E21").Value = "['ELM-2NA-CAD', 'ELM-2NA-CAD']"
